$p = $ppt.ActivePresentation
$hm = $p.HandoutMaster
$dateShape = $hm.Shapes.Item(2)
Write-Host "Type: $($dateShape.Type)"
Write-Host "HasTextFrame: $($dateShape.HasTextFrame)"
$tf = $dateShape.TextFrame
Write-Host "HasText: $($tf.HasText)"
$tr = $tf.TextRange
Write-Host "Before Text=[$($tr.Text)] Length=$($tr.Length)"
$tr.Text = "18.05.2021"
Write-Host "Immediately after set, tr.Text=[$($tr.Text)]"
$tr2 = $dateShape.TextFrame.TextRange
Write-Host "Re-fetched tr2.Text=[$($tr2.Text)]"
$tr3 = $p.HandoutMaster.Shapes.Item(2).TextFrame.TextRange
Write-Host "Re-fetched via hm again tr3.Text=[$($tr3.Text)]"
